# news and keywords.xlsx -- wip political leanings + test run of analysis
# on test data.
#
# Inserts a new "Score (v2)" column (C) between "News Sites (v2)" (B) and
# "Keywords" (was C, now D) -- shifting the former columns C:G to D:H.
# Fills in the rest of the "News Sites (v2)" list (rows 8-31) plus the new
# numeric leaning scores in the inserted column, and appends one more
# keyword ("hamas") to the bottom of "Keywords (v2)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; everything from the old C onward (the
# Keywords* columns) slides right by one.
$ws.Columns("C:C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Score (v2)"

# --- Column B: "News Sites (v2)" -- fill in the remaining site slugs ---
$newsSitesV2 = @{
    8  = "abcnews.go"
    9  = "cbsnews"
    10 = "foxnews"
    11 = "nytimes"
    12 = "usatoday"
    13 = "vice"
    14 = "politico"
    15 = "npr"
    16 = "buzzfeednews"
    17 = "businessinsider"
    18 = "bloomberg"
    19 = "reuters"
    20 = "huffpost"
    21 = "forbes"
    22 = "msn"
    23 = "nypost"
    24 = "thehill"
    25 = "vox"
    26 = "bbc"
    27 = "aljazeera"
    28 = "cnbc"
    29 = "theguardian"
    30 = "news.sky"
    31 = "time"
}
foreach ($row in $newsSitesV2.Keys) {
    $ws.Cells.Item($row, 2).Value = $newsSitesV2[$row]
}

# --- Column C: "Score (v2)" -- political-leaning score per site ---
$scoresV2 = @{
    2  = -1
    3  = -1
    4  = 1
    5  = -1
    6  = -1
    7  = -1
    8  = -1
    9  = -1
    10 = 1
    14 = -1
    18 = -1
    19 = 0
    20 = -1
    21 = 0
    23 = 1
    25 = -1
    26 = 0
    27 = -1
    28 = -1
    29 = -1
    31 = -1
}
foreach ($row in $scoresV2.Keys) {
    $ws.Cells.Item($row, 3).Value = $scoresV2[$row]
}

# --- Column E (was D): "Keywords (v2)" -- one more keyword appended ---
$ws.Range("E24").Value = "hamas"
